# Fix mojibake "Â±" (U+00C2 U+00B1) -> "±" (U+00B1) in columns B, C, D for rows 2-17.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$badChar = [string][char]194 + [string][char]177
$goodChar = [string][char]177

for ($row = 2; $row -le 17; $row++) {
    foreach ($col in @("B", "C", "D")) {
        $cell = $ws.Range("$col$row")
        $val = $cell.Value2
        if ($val -ne $null -and $val.Contains($badChar)) {
            $cell.Value2 = $val.Replace($badChar, $goodChar)
        }
    }
}
